# Weekly data refresh: a new price-survey row for Brócoli at Terminal
# Hortofrutícola Agro Chillán is inserted at row 91 (pushing the existing
# rows 91-148 down to 92-149), matching the "Fruta / hortaliza, semanal"
# weekly update pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 91; existing rows 91..148 shift down
# to 92..149 and the sheet's used-range dimension grows to A1:R149
# automatically.
$ws.Rows(91).Insert()

# Populate the newly inserted row 91 with this week's record.
$ws.Range("A91").Value = 7
$ws.Range("B91").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C91").Value = 'Ñuble'
$ws.Range("D91").Value = 44438
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112023
$ws.Range("G91").Value = 'Brócoli'
$ws.Range("H91").Value = 'Sin especificar'
$ws.Range("I91").Value = 'Primera'
$ws.Range("J91").Value = 300
$ws.Range("K91").Value = 700
$ws.Range("L91").Value = 750
$ws.Range("M91").Value = 725
$ws.Range("N91").Value = '$/unidad'
$ws.Range("O91").Value = 'Región del Maule'
$ws.Range("P91").Value = 725
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = 'Hortaliza'
